$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32 (pushes the existing rows 32-40 down to 33-41)
$ws.Rows(32).Insert()

# Populate the new product row (category 5, next sequence number 6)
$ws.Range("A32").Value = "20093182"
$ws.Range("B32").Value = "IDM CANDY GUMMY 100"
$ws.Range("C32").Value = "SHB01D"
$ws.Range("D32").Value = "5"
$ws.Range("E32").Value = "6"
$ws.Range("F32").Value = "PT,(E-1B)"

# Copy the formatting (borders etc.) from the row above onto the new row
$ws.Range("A31:F31").Copy()
$ws.Range("A32:F32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
